# Update on slides & schedule
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Presentation Slides task is now in-charge of "All" (was "Wesley, Russell, Boon Jun")
$ws.Range("D16").Value = "All"

# Readme task is now in-charge of "All" (was "Sathya")
$ws.Range("D17").Value = "All"

# Update the active selection to match the saved cursor position
$ws.Range("F13").Select()
